$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Insert 4 blank rows right below row 39 (they become rows 40-43; the old
#    rows 40+ shift down to 44+).
# ---------------------------------------------------------------------------
$ws.Rows("40:43").Insert()

# ---------------------------------------------------------------------------
# 2) Propagate row 39's formatting (styles) into the 4 new rows, column by
#    column, so the inserted cells keep the same look (borders/number format)
#    as the rest of the results table instead of Excel's blank default style.
# ---------------------------------------------------------------------------
$resultCols = "A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T"
foreach ($col in $resultCols) {
    $ws.Range("$col" + "39").Copy()
    $ws.Range("$col" + "40:" + "$col" + "43").PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Row 39 already existed as "match 30" but had no result recorded yet
#    (format 1 / blank scores). The commit records its result and relabels
#    it "CSK vs MI" (format 2). Also repair P39, whose VLOOKUP formula was
#    missing in the source file (every other row in the block has it).
# ---------------------------------------------------------------------------
$ws.Range("P39").Formula = "=IF(ISERROR(VLOOKUP(RANK(Q39, (`$T39,`$Q39,`$N39,`$K39,`$H39,`$E39), 0),  `$A`$2:`$C`$7, `$B39+1, FALSE)),"""",VLOOKUP(RANK(Q39, (`$T39,`$Q39,`$N39,`$K39,`$H39,`$E39), 0),  `$A`$2:`$C`$7, `$B39+1, FALSE))"

$ws.Range("B39").Value = 2
$ws.Range("C39").Value = "CSK vs MI"
$ws.Range("E39").Value = 80
$ws.Range("H39").Value = 20
$ws.Range("K39").Value = 60
$ws.Range("N39").Value = 100
$ws.Range("Q39").Value = 0
$ws.Range("T39").Value = 40

# ---------------------------------------------------------------------------
# 4) Fill in the newly inserted rows 40-43 as the next four (still unplayed)
#    matches: match 31 "KKR vs RCB", 32 "PBKS vs RR", 33 "DC vs SRH",
#    34 "MI vs KKR" - same VLOOKUP formulas as the rest of the table, no
#    scores entered yet.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 40; Match = 31; Name = "KKR vs RCB" },
    @{ Row = 41; Match = 32; Name = "PBKS vs RR" },
    @{ Row = 42; Match = 33; Name = "DC vs SRH" },
    @{ Row = 43; Match = 34; Name = "MI vs KKR" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Match
    $ws.Range("B$r").Value = 2
    $ws.Range("C$r").Value = $item.Name

    $ws.Range("D$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(E$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(E$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
    $ws.Range("G$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(H$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(H$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
    $ws.Range("J$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(K$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(K$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
    $ws.Range("M$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(N$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(N$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
    $ws.Range("P$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(Q$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(Q$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
    $ws.Range("S$r").Formula = "=IF(ISERROR(VLOOKUP(RANK(T$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE)),"""",VLOOKUP(RANK(T$r, (`$T$r,`$Q$r,`$N$r,`$K$r,`$H$r,`$E$r), 0),  `$A`$2:`$C`$7, `$B$r+1, FALSE))"
}

# ---------------------------------------------------------------------------
# 5) The summary block (Prize/Total rows, old rows 40-45) has already shifted
#    down automatically to rows 44-49 by the row insert. Its SUM ranges still
#    read D10:D39 etc. (Excel does not extend a SUM when rows are appended
#    right after the range's last row) - widen them to D10:D43 to cover the
#    four freshly added match rows, now that the block lives at row 46.
# ---------------------------------------------------------------------------
$totalCols = "E","H","K","N","Q","T"
$srcCols   = "D","G","J","M","P","S"
for ($i = 0; $i -lt $totalCols.Count; $i++) {
    $ws.Range($totalCols[$i] + "46").Formula = "=SUM(" + $srcCols[$i] + "10:" + $srcCols[$i] + "43)"
}

# ---------------------------------------------------------------------------
# 6) Re-point the six conditional-formatting rules (row 42 -> row 46) so they
#    keep tracking the Total cells without losing their original dxf ids.
# ---------------------------------------------------------------------------
foreach ($col in $totalCols) {
    $oldRange = $ws.Range($col + "42")
    $newRange = $ws.Range($col + "46")
    $fcs = $oldRange.FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------------
# 7) Cosmetic: move the active-cell selection to the new Total/U cell, same
#    as the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("U46").Select()

Write-Output "edit complete"
